$wb = $excel.ActiveWorkbook

# ---------------- Sheet "Classi": new row 12 ----------------
$ws1 = $wb.Worksheets.Item("Classi")
$ws1.Range("A12").Value = "21  - 26 marzo 2022"
$ws1.Range("B12").Value = 5488
$ws1.Range("B12").NumberFormat = "#,##0"
$ws1.Range("C12").Value = 8157
$ws1.Range("C12").NumberFormat = "#,##0"
$ws1.Range("D12").Value = 0.6720000000000000417443857
$ws1.Range("D12").NumberFormat = "0.0%"
$ws1.Range("E12").Value = 376516
$ws1.Range("E12").NumberFormat = "#,##0"
$ws1.Range("F12").Value = 253348
$ws1.Range("F12").NumberFormat = "#,##0"
$ws1.Range("G12").Value = 0.6729999999999999316102617
$ws1.Range("G12").NumberFormat = "0.0%"
$ws1.Range("H12").Value = 253112
$ws1.Range("H12").NumberFormat = "#,##0"
$ws1.Range("I12").Value = 24396
$ws1.Range("I12").NumberFormat = "#,##0"
$ws1.Range("J12").Value = 0.9990000000000001101341240
$ws1.Range("J12").NumberFormat = "0.0%"
$ws1.Range("K12").Value = 0.0960000000000000019984014
$ws1.Range("K12").NumberFormat = "0.0%"
$ws1.Range("N12").Value = 236
$ws1.Range("O12").Value = 0.0010000000000000000208167
$ws1.Range("O12").NumberFormat = "0.0%"
$ws1.Range("O12").Font.Color = 0
$null = $ws1.Range("O12").Select()

# ---------------- Sheet "Alunni in presenza": fix F11 style + new row 12 ----------------
$ws2 = $wb.Worksheets.Item("Alunni in presenza")
$ws2.Range("F11").NumberFormat = "0.0%"
$ws2.Range("A12").Value = "21  - 26 marzo 2022"
$ws2.Range("B12").Value = 7393525
$ws2.Range("B12").NumberFormat = "#,##0"
$ws2.Range("C12").Value = 4960243
$ws2.Range("C12").NumberFormat = "#,##0"
$ws2.Range("D12").Value = 0.6709999999999999298339048
$ws2.Range("D12").NumberFormat = "0.0%"
$ws2.Range("E12").Value = 4786890
$ws2.Range("E12").NumberFormat = "#,##0"
$ws2.Range("F12").Value = 0.9649999999999999689137553
$ws2.Range("F12").NumberFormat = "0.0%"
# selection on this sheet is unchanged from before (F11) per target diff

# ---------------- Sheet "Alunni": new rows 32-34 ----------------
$ws3 = $wb.Worksheets.Item("Alunni")
$ws3.Range("A32").Value = "21  - 26 marzo 2022"
$ws3.Range("B32").Value = "Infanzia"
$ws3.Range("C32").Value = 0.9649999999999999689137553
$ws3.Range("C32").NumberFormat = "0.0%"
$ws3.Range("D32").Value = 5516357
$ws3.Range("D32").NumberFormat = "#,##0"
$ws3.Range("E32").Value = 14530
$ws3.Range("E32").NumberFormat = "#,##0"
$ws3.Range("F32").Value = 0.0260000000000000022759572
$ws3.Range("F32").NumberFormat = "0.0%"
$ws3.Range("A33").Value = "21  - 26 marzo 2022"
$ws3.Range("B33").Value = "Primaria"
$ws3.Range("C33").Value = 1566311
$ws3.Range("C33").NumberFormat = "#,##0"
$ws3.Range("D33").Value = 1506900
$ws3.Range("D33").NumberFormat = "#,##0"
$ws3.Range("E33").Value = 59411
$ws3.Range("E33").NumberFormat = "#,##0"
$ws3.Range("F33").Value = 0.0379999999999999990563104
$ws3.Range("F33").NumberFormat = "0.0%"
$ws3.Range("A34").Value = "21  - 26 marzo 2022"
$ws3.Range("B34").Value = "Sec. 1° e 2° Grado"
$ws3.Range("C34").Value = 2827767
$ws3.Range("C34").NumberFormat = "#,##0"
$ws3.Range("D34").Value = 2728355
$ws3.Range("D34").NumberFormat = "#,##0"
$ws3.Range("E34").Value = 99412
$ws3.Range("E34").NumberFormat = "#,##0"
$ws3.Range("F34").Value = 0.0350000000000000033306691
$ws3.Range("F34").NumberFormat = "0.0%"
$null = $ws3.Range("F35").Select()

# ---------------- Sheet "Personale scolastico": new row 12 ----------------
$ws4 = $wb.Worksheets.Item("Personale scolastico")
$ws4.Range("A12").Value = "21  - 26 marzo 2022"
$ws4.Range("B12").Value = 775867
$ws4.Range("B12").NumberFormat = "#,##0"
$ws4.Range("C12").Value = 517395
$ws4.Range("C12").NumberFormat = "#,##0"
$ws4.Range("D12").Value = 0.6670000000000000373034936
$ws4.Range("D12").NumberFormat = "0.0%"
$ws4.Range("E12").Value = 491320
$ws4.Range("E12").NumberFormat = "#,##0"
$ws4.Range("F12").Value = 0.9499999999999999555910790
$ws4.Range("F12").NumberFormat = "0.0%"
$ws4.Range("G12").Value = 204526
$ws4.Range("G12").NumberFormat = "#,##0"
$ws4.Range("H12").Value = 137225
$ws4.Range("H12").NumberFormat = "#,##0"
$ws4.Range("I12").Value = 0.6709999999999999298339048
$ws4.Range("I12").NumberFormat = "0.0%"
$ws4.Range("J12").Value = 132048
$ws4.Range("J12").NumberFormat = "#,##0"
$ws4.Range("K12").Value = 0.9620000000000000772715225
$ws4.Range("K12").NumberFormat = "0.0%"
$null = $ws4.Range("K13").Select()

